{"js": "// Toll receipt template: replace the hardcoded Peaje value \"$11000\"\n// with the template expression \"$${value}\" (renders as a literal \"$\"\n// immediately followed by the \"${value}\" placeholder), matching the\n// pattern already used elsewhere in the document (${code}, ${date},\n// ${time}).\n//\n// \"Logica agregada para modificar el valor del Peaje\"\n\nconst body = context.document.body;\n\n// Search for the literal \"$11000\" only (not the surrounding \"Valor: \"\n// label), so the existing \"Valor:\" run/formatting is left untouched\n// and only the value itself is replaced.\nconst results = body.search(\"$11000\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"$11000\" in the document body.');\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  const found = results.items[i];\n  // \"$\" + \"${value}\" written this way (string concatenation, not\n  // String.replace) so the literal \"$$\" is preserved as-is instead of\n  // being collapsed/interpreted as a replacement-pattern escape.\n  found.insertText(\"$\" + \"${value}\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Toll receipt template: replace the hardcoded Peaje value \"$11000\"\n# with the template expression \"$${value}\" (renders as a literal \"$\"\n# immediately followed by the \"${value}\" placeholder), matching the\n# pattern already used elsewhere in the document (${code}, ${date},\n# ${time}).\n#\n# \"Logica agregada para modificar el valor del Peaje\"\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"`$11000\"\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0          # wdFindStop - don't loop back, so the search below terminates\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$rng.Find.MatchWildcards = $false\n\n$foundAny = $false\n$guard = 0\n\n# Replace every occurrence of the literal \"$11000\" with \"$${value}\".\n# Only the matched value text is touched; the surrounding \"Valor: \"\n# label and trailing space keep their own runs/formatting untouched.\nwhile ($rng.Find.Execute() -and $guard -lt 25) {\n    $foundAny = $true\n    $guard = $guard + 1\n    $rng.Text = \"`$`${value}\"\n    # Continue searching after the just-replaced text.\n    $rng.Collapse(0)  # wdCollapseEnd\n    $rng.End = $d.Content.End\n}\n\nif (-not $foundAny) {\n    throw 'Could not find \"$11000\" in the document.'\n}\n"}
